$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reduce stored numeric precision ("custom accuracy") on row 5 ---
# Values are rounded to 2 decimal places.
$ws.Range("B5").Value  = 13.63
$ws.Range("C5").Value  = 10.29
$ws.Range("D5").Value  = 0.4
$ws.Range("E5").Value  = 28.82
$ws.Range("F5").Value  = 24.06
$ws.Range("G5").Value  = 10.24
$ws.Range("H5").Value  = 36.71
$ws.Range("I5").Value  = 16.08
$ws.Range("J5").Value  = 7.61
$ws.Range("K5").Value  = 10.9
$ws.Range("L5").Value  = 11.76
$ws.Range("M5").Value  = 12.62
$ws.Range("N5").Value  = 3.69
$ws.Range("O5").Value  = 10.35
$ws.Range("P5").Value  = 14.99
$ws.Range("Q5").Value  = 8.44
$ws.Range("R5").Value  = 0.37
$ws.Range("S5").Value  = 0.36
$ws.Range("T5").Value  = 152.68
$ws.Range("U5").Value  = 28.92
$ws.Range("V5").Value  = 9.92
$ws.Range("W5").Value  = 19.71
$ws.Range("X5").Value  = 10.21
$ws.Range("Y5").Value  = 1.38
$ws.Range("Z5").Value  = 18.32
$ws.Range("AA5").Value = 8.529999999999999
$ws.Range("AB5").Value = 7.58
$ws.Range("AC5").Value = 9.4
$ws.Range("AD5").Value = 12.37
$ws.Range("AE5").Value = 0.57
$ws.Range("AF5").Value = 32.89
$ws.Range("AG5").Value = 5.4
$ws.Range("AH5").Value = 11.98

# --- Remove the now-superfluous last data row (row 6) ---
$ws.Rows.Item(6).Delete()

# --- Narrow a few columns by one character unit (8 -> 7) ---
# Excel's ColumnWidth property pads by 5/6 of a character width relative to
# the raw stored column width, so subtract 5/6 from the desired stored width.
$ws.Columns.Item(5).ColumnWidth  = 7 - 0.8333333333333333
$ws.Columns.Item(7).ColumnWidth  = 7 - 0.8333333333333333
$ws.Columns.Item(24).ColumnWidth = 7 - 0.8333333333333333
